$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (A2 = 0)
$ws.Range("B2").Value = 9.885135962337825
$ws.Range("C2").Value = 8.617752600969302
$ws.Range("D2").Value = 10.33143707665598
$ws.Range("E2").Value = 9.026718694503387
$ws.Range("F2").Value = 10.48677899168146
$ws.Range("G2").Value = 11.13231966990435
$ws.Range("H2").Value = 10.24089607922481
$ws.Range("I2").Value = 8.877942761470294
$ws.Range("J2").Value = 10.83344898740094
$ws.Range("K2").Value = 10.94868145524048
$ws.Range("L2").Value = 8.747939589530533
$ws.Range("M2").Value = 9.663079324017239
$ws.Range("N2").Value = 12.0943182927827
$ws.Range("O2").Value = 8.9015775901431
$ws.Range("P2").Value = 9.402516490669564

# Row 3 (A3 = 1)
$ws.Range("B3").Value = -0.3933032992204761
$ws.Range("C3").Value = -0.1466215587787953
$ws.Range("D3").Value = -0.4509968056791822
$ws.Range("E3").Value = -0.2216401644031936
$ws.Range("F3").Value = -0.4808882761941704
$ws.Range("G3").Value = -0.5601272684600538
$ws.Range("H3").Value = -0.4137386807976613
$ws.Range("I3").Value = -0.2738042408284731
$ws.Range("J3").Value = -0.5377964605821706
$ws.Range("K3").Value = -0.520459718111276
$ws.Range("L3").Value = -0.1821984401405221
$ws.Range("M3").Value = -0.3405392501152689
$ws.Range("N3").Value = -0.7308316978530964
$ws.Range("O3").Value = -0.2692360246771852
$ws.Range("P3").Value = -0.2791351663839501

# Row 4 (A4 = 2)
$ws.Range("B4").Value = 0.2756558215827666
$ws.Range("C4").Value = 0.2573378951961205
$ws.Range("D4").Value = 0.2346451430030998
$ws.Range("E4").Value = 0.1257363849853937
$ws.Range("F4").Value = 0.3679740793016473
$ws.Range("G4").Value = 0.5263917086614048
$ws.Range("H4").Value = 0.4652642320614955
$ws.Range("I4").Value = 0.3325394677740252
$ws.Range("J4").Value = 0.5945383818905342
$ws.Range("K4").Value = 0.5621024745560772
$ws.Range("L4").Value = 0.2567723425716453
$ws.Range("M4").Value = 0.3962517170067417
$ws.Range("N4").Value = 0.3990141238497889
$ws.Range("O4").Value = 0.302725537219417
$ws.Range("P4").Value = 0.8482756995338175

# Row 5 (A5 = 3)
$ws.Range("B5").Value = 0.04424702496362989
$ws.Range("C5").Value = 0.0413067143202412
$ws.Range("D5").Value = 0.03766417643726536
$ws.Range("E5").Value = 0.02018263548123487
$ws.Range("F5").Value = 0.05906553389419372
$ws.Range("G5").Value = 0.08449401481911271
$ws.Range("H5").Value = 0.07468210891576581
$ws.Range("I5").Value = 0.05337773041579533
$ws.Range("J5").Value = 0.09543261039908044
$ws.Range("K5").Value = 0.09022614534673694
$ws.Range("L5").Value = 0.04121593437244355
$ws.Range("M5").Value = 0.06360453232443045
$ws.Range("N5").Value = 0.06404794136923928
$ws.Range("O5").Value = 0.04859213321005088
$ws.Range("P5").Value = 0.1361613763054296

# Row 7: significance stars update
$ws.Range("M7").Value = "***"
